# Add the new dialog system.
#
# phoneMessageData.xlsx / Sheet1 stores a small "data table" describing
# phone messages (var name / type / Chinese label / sample rows). This
# change drops the unused "iconpath" ("图标的路径") column and introduces a
# new "condition" ("条件") column in its place at the right edge of the
# table, then moves the active selection onto the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D was "iconpath" / "图标的路径". Deleting it shifts the old
# "content" column (E) and the following (unused) column (F) one slot to
# the left, landing on D and E respectively - exactly the column layout
# the new table needs (B/C id+sender columns are untouched).
$ws.Columns.Item(4).Delete()

# The freed-up column E becomes the new "condition" field: a ##var header,
# a ##type of "string" (matching every other field), and the Chinese label
# "条件" on the description row. The two existing message rows (4 and 5)
# are left blank for this new column - no condition data for them yet.
$ws.Range("E1").Value = "condition"
$ws.Range("E2").Value = "string"
$ws.Range("E3").Value = "条件"

# Park the selection on the newly added header cell.
$ws.Range("E3").Select()
